# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Mapping of row -> new F value (old value in comment for reference).
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 294   # was 291
    4  = 56    # was 55
    7  = 2075  # was 2068
    8  = 72    # was 71
    10 = 4535  # was 4520
    12 = 286   # was 285
    13 = 108   # was 107
    15 = 134   # was 133
    18 = 86    # was 84
    19 = 3415  # was 3391
    20 = 84    # was 81
    21 = 539   # was 533
    22 = 21    # was 20
    24 = 86    # was 84
    25 = 95    # was 94
    28 = 65    # was 64
    29 = 207   # was 206
    31 = 668   # was 654
    32 = 2066  # was 2046
    33 = 391   # was 386
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
